$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows starting at row 12 (pushes old rows 12-22 down to 15-25)
$ws.Range("A12:A14").EntireRow.Insert()

# Row 12: new label-only row "Docentes responsáveis:" -- copy style/formatting from a
# similar label-only row (row 19 after the shift -> "Syllabus:" pattern has style s=1, no row height)
$ws.Cells.Item(19, 1).Copy($ws.Cells.Item(12, 1))
$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"

# Row 13: B/C pair "5983729 - Fernando Vernilli Junior" -- copy style from row 10's B/C
# (s=2 / s=3, no custom row height)
$ws.Cells.Item(10, 2).Copy($ws.Cells.Item(13, 2))
$ws.Cells.Item(10, 3).Copy($ws.Cells.Item(13, 3))
$ws.Cells.Item(13, 2).Value = "5983729 - Fernando Vernilli Junior"
$ws.Cells.Item(13, 3).Value = "5983729 - Fernando Vernilli Junior"

# Row 14: B/C pair "1922320 - Sebastiao Ribeiro"
$ws.Cells.Item(10, 2).Copy($ws.Cells.Item(14, 2))
$ws.Cells.Item(10, 3).Copy($ws.Cells.Item(14, 3))
$ws.Cells.Item(14, 2).Value = "1922320 - Sebastiao Ribeiro"
$ws.Cells.Item(14, 3).Value = "1922320 - Sebastiao Ribeiro"
